$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    current "2022-Q2" sheet) holding the fund-detail breakdown for the new
#    quarter. Clone the layout/formatting of the existing "2022-Q2" sheet
#    (same headers, same two funds, same bold/bordered header & index
#    styling) and then overwrite the quarter-specific numbers.
# ---------------------------------------------------------------------------
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")
$totalSheet = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($q2Sheet, $totalSheet)
$newSheet.Name = "2022-Q3"

# Re-resolve both sheets by name: inserting a sheet between them shifts
# internal positions, so any reference captured before the Add() call can
# end up stale.
$q2Sheet  = $wb.Worksheets.Item("2022-Q2")
$q3Sheet  = $wb.Worksheets.Item("2022-Q3")

$q2Sheet.Range("A1:H3").Copy($q3Sheet.Range("A1"))

# Fund code / name / size / position columns are stored as text - keep them
# text after overwriting so leading zeros and fixed-decimal strings survive.
$q3Sheet.Range("D2:G3").NumberFormat = "@"

# Row 2 - 009837 华夏磐锐一年定期开放混合A
$q3Sheet.Range("D2").Value = "14.02"
$q3Sheet.Range("E2").Value = "94.15"
$q3Sheet.Range("F2").Value = "4.42"
$q3Sheet.Range("G2").Value = "0.6197"

# Row 3 - 009838 华夏磐锐一年定期开放混合C
$q3Sheet.Range("D3").Value = "0.39"
$q3Sheet.Range("E3").Value = "94.15"
$q3Sheet.Range("F3").Value = "4.42"
$q3Sheet.Range("G3").Value = "0.0172"

# ---------------------------------------------------------------------------
# 2. Insert a new row into "总计" for the 2022-Q3 summary entry, shifting the
#    existing quarters down by one row, then renumber the index column (A)
#    for every data row (it's a plain 0-based row counter, not business
#    data, so it has to stay 0,1,2,3 after the insert).
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.64

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
